$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add new row 16 with the test-sync data for Testmail #14
$ws.Range("A16").Value = "Heb je de CE-certificaten van dit product?"
$ws.Range("B16").Value = "mailmind.test@zohomail.eu"
$ws.Range("C16").Value = "Testmail #14: Heb je de CE-certificaten van dit product?"
$ws.Range("D16").Value = "Productinformatie"
$ws.Range("E16").Value = "Beste sender,`nBedankt voor je e-mail. Om je vraag over de CE-certificaten van het product te beantwoorden, heb ik wat meer informatie nodig. Kun je het productnummer of de naam van het product doorgeven, zodat ik de relevante certificaten kan opzoeken? `nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Range("F16").Value = "2025-07-27 19:41:20"
$ws.Range("G16").Value = "Ja"
$ws.Range("H16").Value = "Nee"
$ws.Range("I16").Value = "Ja"
$ws.Range("J16").Value = "Nee"

# Setting a multi-line value auto-expands the row height (Excel's own
# wrap-based autofit); restore the default so row 16 matches the plain,
# unattributed <row> of the other entries.
$ws.Rows.Item(16).EntireRow.AutoFit()

# Extend the conditional-formatting ranges so the new row inherits the
# same colour rules as the rest of the log (D/G/H/I/J now go to row 16).
$cols = "D", "G", "H", "I", "J"
foreach ($col in $cols) {
    $oldRange = $ws.Range($col + "2:" + $col + "15")
    $newRange = $ws.Range($col + "2:" + $col + "16")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Keep the Dashboard summary in sync: one more "Productinformatie" mail.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 5
